# Add a second worksheet "Лист1" with correlation data + a clustered
# column chart, matching the target diff.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new worksheet, placed after "List1" ------------------------
$ws1 = $wb.Worksheets.Item("List1")
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Лист1"

# --- 2. Header row (bold, left aligned) -------------------------------------
# NB: shared-string table order must match the original authoring order
# (Code Reuse registered before Development Mode), so columns are written
# out of left-to-right sequence on purpose: A,B,C,E,D,F,G,H.
$headerCols = @(1, 2, 3, 5, 4, 6, 7, 8)
$headerText = @("Language", "Team", "Size", "Code Reuse", "Development Mode", "Architecture", "Customer Quality", "PM Quality")
for ($i = 0; $i -lt $headerCols.Count; $i++) {
    $cell = $ws2.Cells.Item(1, $headerCols[$i])
    $cell.Value = $headerText[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4131
}

# --- 3. Data row (values, 2-decimal number format, left aligned) -----------
$values = @(-0.22540381906681656, 0.58833261402197501, 0.88199951980579872, -0.26789662524016455, -0.42312093918262578, 0.20775771477150626, -0.2546595209141273, -0.2971228931536446)
for ($i = 0; $i -lt $values.Count; $i++) {
    $cell = $ws2.Cells.Item(2, $i + 1)
    $cell.Value = $values[$i]
    $cell.NumberFormat = "0.00"
    $cell.HorizontalAlignment = -4131
}

# --- 4. Chart ----------------------------------------------------------------
$chartObj = $ws2.Shapes.AddChart2(201, 51).Chart
$chartObj.SetSourceData($ws2.Range("A1:H2"))
$chartObj.ChartType = 51
$chartObj.HasTitle = $false

$chartObj.Axes(1).AxisTitle.Text = "Feature"
$chartObj.Axes(2).HasTitle = $true
$chartObj.Axes(2).AxisTitle.Text = "Correlation ratio"
$chartObj.SeriesCollection(1).Name = "Correlation"

# --- 5. Column widths on the new sheet (auto-fit-like, from typed headers) -
$ws2.Columns.Item(4).ColumnWidth = 17.0
$ws2.Columns.Item(5).ColumnWidth = 17.666666666666668
$ws2.Columns.Item(7).ColumnWidth = 18.5
$ws2.Columns.Item(8).ColumnWidth = 15.499999999999998

# --- 6. Selections -------------------------------------------------------------
$ws1.Range("H1:L1").Select()
$ws2.Range("B9").Select()
$ws1.Activate()
